$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers corresponding to B, C, D, E, F, H, J, L, N
$cols = @(2,3,4,5,6,8,10,12,14)

# New values for rows 2-25 (data rows 0-23), one row per array, values in column order above
$data = @(
    @(4.576389969012894,0.3342714604671073,0.006836383813787705,0.05039595700296362,4.735664047197929,0.07973214163530429,0.1416895320636762,0.4414515150750873,2.757242386420714),
    @(4.455837958333063,0.3080797494908154,0.006087666367751154,0.05047652714795842,4.707275549390317,0.07973214163530429,0.1424044873738985,0.4365187045588641,2.774832669723565),
    @(4.384589064129784,0.2922051183744259,0.005626340931364382,0.05053027515289243,4.692248519446323,0.07973214163530429,0.142877273063819,0.4337185642581147,2.786423068807963),
    @(4.356248667022896,0.2857874830918377,0.005437879224803766,0.05055325506179437,4.686727119348774,0.07973214163530429,0.1430784426895304,0.4326349052862213,2.791344288055612),
    @(4.351584614418471,0.284724924407243,0.005406555212072561,0.05055713596217926,4.685846605860064,0.07973214163530429,0.1431123606591438,0.4324584307092039,2.792173396664928),
    @(4.384204048535707,0.2921183606677005,0.005623801248386684,0.05053058070441441,4.692171620306425,0.07973214163530429,0.142879951659916,0.4337037172446827,2.786488637009114),
    @(4.534246999201969,0.3251971288773063,0.006578521943524152,0.05042285066560309,4.72537548266817,0.07973214163530429,0.1419290371534707,0.4397031855138493,2.763143219144894),
    @(4.850603414006684,0.3917435685592068,0.008440870247973464,0.0502454699341871,4.80966801259396,0.07973214163530429,0.1403322143127514,0.4532866387676222,2.723656710721002),
    @(5.096745335844275,0.4417131830359722,0.009807277632454259,0.05013572106776887,4.883454828677458,0.07973214163530429,0.1393219986147578,0.4643832081064403,2.698519053072573),
    @(5.211752386229477,0.4646919755899148,0.01042939837254409,0.05009024535380668,4.919634178634254,0.07973214163530429,0.1388977419359989,0.4696758015430191,2.687932273078744),
    @(5.255742883157154,0.4734299227247902,0.01066512940528952,0.0500736636344089,4.933713011208795,0.07973214163530429,0.1387421577390313,0.4717152914156912,2.684046041560904),
    @(5.246249128182285,0.4715464185729843,0.01061435282271361,0.05007720640071867,4.930664005281756,0.07973214163530429,0.138775440016996,0.4712744789961789,2.68487753961621),
    @(5.215362677103144,0.4654101180308885,0.01044878874493804,0.05008886836531257,4.920784851331206,0.07973214163530429,0.1388848402836693,0.4698428834907844,2.687610086029792),
    @(5.196501210637052,0.4616562183532551,0.01034739727194989,0.05009609483972444,4.914782957005826,0.07973214163530429,0.1389525116279238,0.4689705903443411,2.689299857872143),
    @(5.089290721180021,0.4402165110675469,0.009766637349443386,0.05013878247102599,4.881143240022993,0.07973214163530429,0.1393504347729149,0.4640422581819337,2.6992280606501),
    @(5.024300410456817,0.4271278696604099,0.009410546862596902,0.05016610893299311,4.861177609923999,0.07973214163530429,0.1396035856467783,0.4610816387911996,2.705536529286348),
    @(4.987205250342242,0.4196228631921031,0.009205779138426351,0.05018224525715109,4.849939724228847,0.07973214163530429,0.1397525136695226,0.4594017949646059,2.709244810053463),
    @(4.97469441148894,0.4170857630961109,0.009136454092214308,0.05018778070869473,4.846176907089671,0.07973214163530429,0.1398035088440253,0.4588369801840173,2.710514056279024),
    @(5.031189163502233,0.4285187674723829,0.009448448024095057,0.05016315664026538,4.863277524870512,0.07973214163530429,0.1395762934765941,0.4613944178771732,2.704856715986622),
    @(5.224422813265733,0.4672115044745624,0.01049741436138873,0.05008542562964868,4.923676304184539,0.07973214163530429,0.1388525691257989,0.4702624188602869,2.686804133453364),
    @(5.353277385893648,0.4927118644707775,0.01118386706631469,0.05003834782904892,4.965357868210276,0.07973214163530429,0.1384091392829987,0.4762639787543179,2.675721575942546),
    @(5.284269452310582,0.4790821352789862,0.01081738915385699,0.0500631336759344,4.942908742006523,0.07973214163530429,0.13864310173037,0.473041963680032,2.681570796080734),
    @(5.028073921757198,0.4278898808726694,0.009431313059273805,0.05016449004497581,4.862327404045857,0.07973214163530429,0.1395886217128464,0.4612529411394206,2.705163805911752),
    @(4.762628969091793,0.3735556747814996,0.007937684988892357,0.05028983783157026,4.78479449909409,0.07973214163530429,0.1407355528815408,0.4494164429837326,2.733661364364735)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $rowIndex = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($rowIndex, $cols[$c]).Value = $rowVals[$c]
    }
}

Write-Host "Updated pl_mw values for 380 kV case"
